# Tools refactoring (z0lib): the sample asset-import sheet gains a third
# "purchase_date" data row. F4 was previously blank (using a date-formatted
# style); it now holds the same kind of free-text date string already used
# in F2/F3 ("<2-12-20", "<1-03-31"), so give it the identical look & feel
# before writing the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F4").Value = "<2-10-01"

# Leave the cursor on the newly-filled cell, matching the saved selection.
$ws.Range("F4").Select()
